$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 193; this shifts old rows 193:223 down to 194:224
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new weekly price observation
$ws.Cells.Item(193, 1).Value = 3
$ws.Cells.Item(193, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(193, 3).Value = "Coquimbo"
$ws.Cells.Item(193, 4).Value = 44505
$ws.Cells.Item(193, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(193, 5).Value = 5
$ws.Cells.Item(193, 6).Value = 100112009
$ws.Cells.Item(193, 7).Value = "Acelga"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 270
$ws.Cells.Item(193, 11).Value = 2000
$ws.Cells.Item(193, 12).Value = 2200
$ws.Cells.Item(193, 13).Value = 2104
$ws.Cells.Item(193, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(193, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(193, 16).Value = 351
$ws.Cells.Item(193, 17).Value = 6
$ws.Cells.Item(193, 18).Value = "Hortaliza"
